# Daily attendance processing - 2025-11-13 04:28:23
# Reorders the "Recorded By" (column G) contributor lists for specific
# rows so that `backup@backdoor.com` / rotated entries appear first,
# matching the upstream attendance export ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact literal replacements observed for the "Recorded By" column (G).
# Any cell whose text matches a key below is replaced with its value;
# all other cell contents are left untouched.
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value()
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
